$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 19.92000007629395
$ws.Range("L2").Value = "19.92±0.22"
$ws.Range("O2").Value = 87.31
$ws.Range("C3").Value = 25.28000068664551
$ws.Range("F3").Value = 0.29
$ws.Range("L3").Value = "25.28±0.22"
$ws.Range("O3").Value = 84.93000000000001
$ws.Range("C4").Value = 30.06999969482422
$ws.Range("F4").Value = 0.1
$ws.Range("L4").Value = "30.07±0.1"
$ws.Range("O4").Value = 108.28
$ws.Range("C5").Value = 31.69000053405762
$ws.Range("F5").Value = 0.31
$ws.Range("L5").Value = "31.69±0.12"
$ws.Range("O5").Value = 90.29000000000001
$ws.Range("C6").Value = 34.84999847412109
$ws.Range("F6").Value = 0.15
$ws.Range("L6").Value = "34.85±0.1"
$ws.Range("O6").Value = 105
$ws.Range("C7").Value = 37.13999938964844
$ws.Range("F7").Value = 0.14
$ws.Range("I7").Value = 0.09000000357627869
$ws.Range("L7").Value = "37.14±0.09"
$ws.Range("O7").Value = 107.94
$ws.Range("C8").Value = 39.90000152587891
$ws.Range("F8").Value = 0.11
$ws.Range("L8").Value = "39.9±0.08"
$ws.Range("O8").Value = 111.89
$ws.Range("C9").Value = 42.97999954223633
$ws.Range("L9").Value = "42.98±0.1"
$ws.Range("O9").Value = 114.99
$ws.Range("C10").Value = 45.0099983215332
$ws.Range("L10").Value = "45.01±0.19"
$ws.Range("O10").Value = 106.02
$ws.Range("C11").Value = 46.9900016784668
$ws.Range("F11").Value = 0.16
$ws.Range("L11").Value = "46.99±0.16"
$ws.Range("O11").Value = 110.51
$ws.Range("C12").Value = 53.06999969482422
$ws.Range("F12").Value = 0.17
$ws.Range("L12").Value = "53.07±0.17"
$ws.Range("O12").Value = 108.85
$ws.Range("C13").Value = 55.02999877929688
$ws.Range("F13").Value = 0.14
$ws.Range("I13").Value = 0.1700000017881393
$ws.Range("L13").Value = "55.03±0.15"
$ws.Range("O13").Value = 114.67
$ws.Range("C14").Value = 56.7599983215332
$ws.Range("F14").Value = 0.26
$ws.Range("L14").Value = "56.76±0.21"
$ws.Range("O14").Value = 102.92
$ws.Range("C15").Value = 59.79999923706055
$ws.Range("F15").Value = 0.24
$ws.Range("L15").Value = "59.8±0.18"
$ws.Range("O15").Value = 106.83
$ws.Range("C16").Value = 63.13000106811523
$ws.Range("F16").Value = 0.21
$ws.Range("L16").Value = "63.13±0.17"
$ws.Range("O16").Value = 108.88
$ws.Range("C17").Value = 65.15000152587891
$ws.Range("F17").Value = 0.23
$ws.Range("L17").Value = "65.15±0.22"
$ws.Range("O17").Value = 105.91
$ws.Range("C18").Value = 66.98000335693359
$ws.Range("L18").Value = "66.98±0.2"
$ws.Range("O18").Value = 110.77
$ws.Range("C19").Value = 69.80000305175781
$ws.Range("F19").Value = 0.27
$ws.Range("L19").Value = "69.8±0.22"
$ws.Range("O19").Value = 107.83
$ws.Range("C20").Value = 72.90000152587891
$ws.Range("F20").Value = 0.26
$ws.Range("L20").Value = "72.9±0.22"
$ws.Range("O20").Value = 109.22
